# Append a new row (row 81) of sensor data to each of the four worksheets.
# Columns: A=time, B=总长, C=ID, D=实际长度, E=和校验,
#          F=总长_DEC, G=ID_DEC, H=实际长度_DEC, I=和校验_DEC

$wb = $excel.ActiveWorkbook

$rowsBySheet = @{
    "ROW35-FE-LIFTER"  = @{
        A = "2025-03-07 16:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = "2025-03-07 16:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    }
    "ROW02-FE-LIFTER"  = @{
        A = "2025-03-07 16:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    }
    "ROW02-MID-LIFTER" = @{
        A = "2025-03-07 16:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($rowsBySheet.ContainsKey($name)) {
        $data = $rowsBySheet[$name]
        $newRow = 81

        $ws.Cells.Item($newRow, 1).Value = $data.A
        $ws.Cells.Item($newRow, 2).Value = $data.B
        $ws.Cells.Item($newRow, 3).Value = $data.C
        $ws.Cells.Item($newRow, 4).Value = $data.D
        $ws.Cells.Item($newRow, 5).Value = $data.E
        $ws.Cells.Item($newRow, 6).Value = $data.F

        # G is a very large integer-looking string (25 digits) that exceeds
        # double precision; a plain .Value assignment would be silently
        # coerced to a lossy number. Force text storage, then strip the
        # number-format override so no stray style is left behind.
        $gCell = $ws.Cells.Item($newRow, 7)
        $gCell.NumberFormat = "@"
        $gCell.Value = $data.G
        $gCell.Style = "Normal"

        $ws.Cells.Item($newRow, 8).Value = $data.H
        $ws.Cells.Item($newRow, 9).Value = $data.I
    }
}
